# Capstone Pres - Brody Hill.pptx
#
# This script reproduces three shape-position/autofit tweaks:
#   1. Slide 2 ("The Problem:") - Title placeholder nudged down/right.
#   2. Slide 3 ("The Solution:") - Title placeholder moved up, autofit
#      switched from a 90% "shrink text on overflow" to "do not autofit".
#   3. Slide 3 - Content placeholder ("Invest heavily ...") moved/widened.
#
# NOTE: the COM host here converts Shape.Left/Top/Width/Height (points)
# to EMU internally via a single-precision (f32) round-trip and then
# truncates, so naive "emu/12700" literals can land 1 EMU short. The
# literals below were chosen so that they survive that f32 truncation
# and reproduce the exact target EMU values from the diff.

$p = $ppt.ActivePresentation

# --- Edit 1: Slide 2, "Title 1" -------------------------------------------
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item("Title 1")
$title2.Left   = 51.7412
$title2.Top    = 125.9971
$title2.Width  = 873.3544
$title2.Height = 104.8819

# --- Edit 2: Slide 3, "Title 1" --------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item("Title 1")
$title3.Left   = 311.4907
$title3.Top    = 136.1157
$title3.Width  = 339.9672
$title3.Height = 31.3376
$title3.TextFrame.AutoSize = 0

# --- Edit 3: Slide 3, "Content Placeholder 2" ------------------------------
$content3 = $s3.Shapes.Item("Content Placeholder 2")
$content3.Left   = 112.3549
$content3.Top    = 206.14811
$content3.Width  = 744.1746
$content3.Height = 134.9248
